# ------------------------------------------------------------------
# Adds the "2022-Q4" quarter sheet to the workbook:
#  1. Inserts a new worksheet named "2022-Q4" right after "总计".
#  2. Fills it with the fund-holding detail rows for that quarter.
#  3. Updates the "总计" (summary) sheet by inserting a new row for
#     2022-Q4 at the top of the data and shifting the remaining rows
#     down by one (the table keeps growing historically).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet right after "总计"
# ---------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row (B1:H1), styled like the "总计" header cells.
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Fund detail rows (index, code, name, scale, position, ratio, marketValue, rank)
$q4data = @(
    @(0, "001975", "景顺长城环保优势股票", "40.92", "92.75", "2.96", "1.2112", 10),
    @(1, "010094", "交银施罗德产业机遇混合", "16.65", "85.76", "4.74", "0.7892", 4),
    @(2, "519773", "交银施罗德数据产业灵活配置混合A", "13.79", "86.66", "4.61", "0.6357", 5),
    @(3, "010947", "中欧嘉选混合A", "12.93", "83.49", "2.53", "0.3271", 10),
    @(4, "012703", "华夏核心成长混合A", "4.93", "92.86", "5.01", "0.2470", 10),
    @(5, "501081", "中欧科创主题混合（LOF）A", "7.06", "87.64", "2.93", "0.2069", 9),
    @(6, "910021", "东方红启华三年持有期混合A", "4.18", "91.17", "2.94", "0.1229", 8),
    @(7, "014549", "交银施罗德数据产业灵活配置混合C", "2.43", "86.66", "4.61", "0.1120", 5),
    @(8, "011097", "达诚宜创精选混合A", "0.67", "82.54", "6.50", "0.0436", 3),
    @(9, "010301", "达诚成长先锋混合A", "0.51", "81.22", "6.39", "0.0326", 3),
    @(10, "012710", "华夏核心成长混合C", "0.60", "92.86", "5.01", "0.0301", 10),
    @(11, "011313", "东方红启华三年持有期混合B", "0.89", "91.17", "2.94", "0.0262", 8),
    @(12, "010808", "达诚策略先锋混合A", "0.32", "81.99", "6.44", "0.0206", 3),
    @(13, "010809", "达诚策略先锋混合C", "0.32", "81.99", "6.44", "0.0206", 3),
    @(14, "010302", "达诚成长先锋混合C", "0.30", "81.22", "6.39", "0.0192", 3),
    @(15, "010948", "中欧嘉选混合C", "0.74", "83.49", "2.53", "0.0187", 10),
    @(16, "011031", "达诚价值先锋灵活配置混合C", "0.22", "71.99", "6.17", "0.0136", 2),
    @(17, "011030", "达诚价值先锋灵活配置混合A", "0.20", "71.99", "6.17", "0.0123", 2),
    @(18, "011098", "达诚宜创精选混合C", "0.17", "82.54", "6.50", "0.0110", 3),
    @(19, "017290", "中欧科创主题混合（LOF）C", "0.00", "87.64", "2.93", 0, 9)
)

# Fund code (B) and the scale/position/ratio/value columns (D:G) are
# text-formatted numbers in the source data (e.g. "001975", "40.92"),
# so force those ranges to text before writing, keeping leading zeros
# and the exact formatted precision intact. The very last row's
# "持有市值" (G21) is the lone exception - it is a genuine numeric 0 -
# so it is left out of the text-formatted range.
$q4.Range("B2:B21").NumberFormat = "@"
$q4.Range("D2:G20").NumberFormat = "@"
$q4.Range("D21:F21").NumberFormat = "@"

$rowIdx = 2
foreach ($rec in $q4data) {
    $q4.Range("A$rowIdx").Value = $rec[0]
    $q4.Range("B$rowIdx").Value = $rec[1]
    $q4.Range("C$rowIdx").Value = $rec[2]
    $q4.Range("D$rowIdx").Value = $rec[3]
    $q4.Range("E$rowIdx").Value = $rec[4]
    $q4.Range("F$rowIdx").Value = $rec[5]
    $q4.Range("G$rowIdx").Value = $rec[6]
    $q4.Range("H$rowIdx").Value = $rec[7]
    $rowIdx++
}

# Column A (index numbers) carries the bold/boxed style used
# throughout the workbook for the leading index column.
$total.Range("A2").Copy()
$q4.Range("A2:A21").PasteSpecial(-4122)

[void]$q4.Range("A1").Select()

# ---------------------------------------------------------------
# 2) Update the "总计" sheet: insert the 2022-Q4 summary row and
#    shift the historical rows down by one.
# ---------------------------------------------------------------
$summary = @(
    @(0, "2022-Q4", 20, 3.9),
    @(1, "2022-Q3", 10, 0.41),
    @(2, "2022-Q2", 16, 0.8),
    @(3, "2022-Q1", 17, 3.78),
    @(4, "2021-Q4", 28, 18.08),
    @(5, "2021-Q3", 21, 14.36),
    @(6, "2021-Q2", 8, 1.44),
    @(7, "2021-Q1", 8, 1.69)
)

$rowIdx = 2
foreach ($rec in $summary) {
    $total.Range("A$rowIdx").Value = $rec[0]
    $total.Range("B$rowIdx").Value = $rec[1]
    $total.Range("C$rowIdx").Value = $rec[2]
    $total.Range("D$rowIdx").Value = $rec[3]
    $rowIdx++
}

# Make sure the newly added row 9 carries the same formatting
# (bold index cell) as the rest of column A.
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)
$total.Range("A9").Value = 7

[void]$total.Activate()
[void]$total.Range("A1").Select()

Write-Host "Workbook updated: 2022-Q4 sheet added, summary sheet refreshed."
